$d = $word.ActiveDocument

$replacements = @(
    @{old="549×5=2745"; new="498×6=2988"},
    @{old="694×7=4858"; new="600×9=5400"},
    @{old="217×9=1953"; new="858×7=6006"},
    @{old="249×7=1743"; new="655×3=1965"},
    @{old="511×8=4088"; new="440×2=880"},
    @{old="478×6=2868"; new="430×6=2580"},
    @{old="386×3=1158"; new="754×5=3770"},
    @{old="280×2=560"; new="433×6=2598"},
    @{old="703×4=2812"; new="584×7=4088"},
    @{old="922×4=3688"; new="951×4=3804"},
    @{old="113×9=1017"; new="799×8=6392"},
    @{old="197×5=985"; new="687×5=3435"},
    @{old="422×4=1688"; new="637×4=2548"},
    @{old="949×7=6643"; new="859×3=2577"},
    @{old="264×8=2112"; new="763×5=3815"},
    @{old="783×6=4698"; new="920×2=1840"},
    @{old="381×9=3429"; new="647×9=5823"},
    @{old="648×3=1944"; new="889×7=6223"},
    @{old="779×6=4674"; new="443×9=3987"},
    @{old="812×6=4872"; new="292×9=2628"},
    @{old="813×6=4878"; new="743×2=1486"},
    @{old="270×6=1620"; new="365×2=730"},
    @{old="395×5=1975"; new="426×5=2130"},
    @{old="420×3=1260"; new="489×4=1956"},
    @{old="336×3=1008"; new="815×6=4890"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
